# Adds new fields (processor_full_name, project_name, process_capsule_id)
# to the job upload template by inserting three columns at the front of
# the sheet, filling in their header + sample data, adding a project_name
# dropdown validation, and widening the new columns to match the rest.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 8 columns (platform..modality1.source) three to the
# right, carrying their data, styles, column widths and data validations
# with them (platform -> D, acq_datetime -> E, subject_id -> F,
# metadata_dir -> G, modality0 -> H, modality0.source -> I,
# modality1 -> J, modality1.source -> K).
$ws.Range("A:C").Insert()

# --- New header row (row 1) -------------------------------------------
$ws.Cells.Item(1, 1).Value = "processor_full_name"
$ws.Cells.Item(1, 2).Value = "project_name"
$ws.Cells.Item(1, 3).Value = "process_capsule_id"

# Carry over the bold header formatting (column D already has it, having
# been the original column A) onto the three new header cells.
$ws.Cells.Item(1, 4).Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)

# --- New row 2 data ------------------------------------------------------
$ws.Cells.Item(2, 1).Value = "Anna Apple"
$ws.Cells.Item(2, 2).Value = "Behavior Platform"
$ws.Cells.Item(2, 3).Value = "1f999652-00a0-4c4b-99b5-64c2985ad070"

# --- New row 3 data (no process_capsule_id for this row) -----------------
$ws.Cells.Item(3, 1).Value = "John Smith"
$ws.Cells.Item(3, 2).Value = "Ophys Platform - SLAP2"

# --- New row 4 data (no process_capsule_id for this row) -----------------
$ws.Cells.Item(4, 1).Value = "Anna Apple"
$ws.Cells.Item(4, 2).Value = "Ephys Platform"

# Match the column widths/format used by the rest of the template.
$ws.Range("A:C").ColumnWidth = 12.14

# New dropdown validation for project_name (column B, rows 2-20).
$projectNameList = '"AIND Viral Genetic Tools,Behavior Platform,Brain Computer Interface,Cell Type LUT,Cognitive flexibility in patch foraging,Discovery-Brain Wide Circuit Dynamics,Discovery-Neuromodulator circuit dynamics during foraging,Dynamic Routing,Ephys Platform,Force Foraging,Information seeking in partially observable environments,Learning mFISH/V1omFISH,MSMA Platform,Medulla,Neurobiology of Action,OpenScope,Ophys Platform - FP and indicator testing,Ophys Platform - SLAP2,Single-neuron computations within brain-wide circuits (SCBC),Thalamus in the middle"'
$projectNameRange = $ws.Range("B2:B20")
$projectNameRange.Validation.Add(3, 1, 1, $projectNameList)
$projectNameRange.Validation.IgnoreBlank = $true
$projectNameRange.Validation.InCellDropdown = $true
$projectNameRange.Validation.ShowInput = $true
$projectNameRange.Validation.ShowError = $true
$projectNameRange.Validation.ErrorTitle = "project_name"
$projectNameRange.Validation.ErrorMessage = "Invalid project_name."
$projectNameRange.Validation.InputTitle = "project_name"
$projectNameRange.Validation.InputMessage = "Select a project_name from the dropdown"
